$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("contacts")

# --- Data corrections on the "contacts" sheet ---

# Row 2 (Tintu's record): first name was mistyped, hobbies list had a
# trailing comma that needs to go.
$ws.Range("B2").Value = "rohan"
$ws.Range("L2").Value = "Reading ,Drawing"

# Row 3 (Maya's record): first name + street corrected, hobbies list
# trailing comma removed too.
$ws.Range("B3").Value = "mini"
$ws.Range("H3").Value = "abcd"
$ws.Range("L3").Value = "Reading ,Writing"

# --- Formatting clean-up ---

# The Pincode/Phone header cells (I1, K1) were missing the bordered
# number-style formatting that the data cells below them already use
# (I2:I3, K2:K3) - bring the headers into line with the rest of the
# column by copying that formatting up.
$ws.Range("I2").Copy()
$ws.Range("I1").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("K2").Copy()
$ws.Range("K1").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = 0

# Header/data rows render slightly taller now.
$ws.Rows("1:3").RowHeight = 19.5
